# Fix header typo and replace cluster id data (Imm.* -> MeV.*), dropping the
# stale "leiden_fusion_old2" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column B ("leiden_fusion_old2") is being removed entirely, and the
# old column C ("leiden_fusion") becomes the new column B. Deleting column B
# shifts C into B automatically and keeps the header/style of the surviving
# column intact.
$ws.Columns.Item(2).Delete()

# New set of cluster values for column A and column B (same values in both
# columns), replacing the old Imm.* rows and extending the table down to
# row 28.
$values = @(
    "MeV.2.1",
    "MeV.2.8",
    "MeV.1.4.2",
    "MeV.4.21",
    "MeV.1.4.5",
    "MeV.1.4.7",
    "MeV.1.4.15",
    "MeV.1.4.6",
    "MeV.1.4.4",
    "MeV.1.4.20",
    "MeV.1.4.1",
    "MeV.1.4.11",
    "MeV.1.4.8",
    "MeV.4.12",
    "MeV.4.4",
    "MeV.1.4.0",
    "MeV.3.17",
    "MeV.4.31",
    "MeV.4.1",
    "MeV.4.34",
    "MeV.1.4.13",
    "MeV.3.30",
    "MeV.4.26",
    "MeV.1.4.12",
    "MeV.1.4.21",
    "MeV.4.30",
    "MeV.NA"
)

$row = 2
foreach ($value in $values) {
    $ws.Cells.Item($row, 1).Value = $value
    $ws.Cells.Item($row, 2).Value = $value
    $row = $row + 1
}
